# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# 1) Re-sorts the "Periodo Mora" rows (B16:J62) from descending (2507..2109)
#    to ascending (2109..2507) order.
# 2) Inserts one new trailing data row (row 63) for period 2508, copying the
#    formatting that used to belong to the old "last" row (old row 62), while
#    the former last row (now row 62) takes on the regular interior-row format.
# 3) Updates the summary fields: Valor Mora (E11) and Cant. Periodos (F13).
# 4) The signature block rows shift down by one (old 67/68 -> new 68/69)
#    automatically because of the inserted row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 16
$lastDataRow  = 62

# ---------------------------------------------------------------------------
# Step 1: capture the existing periods (oldest->newest is currently bottom->top)
# and rebuild the ascending list, adding the new period at the end.
# ---------------------------------------------------------------------------
$periods = @()
for ($r = $lastDataRow; $r -ge $firstDataRow; $r--) {
    $periods += $ws.Cells.Item($r, 5).Value()
}
# $periods is now already ascending (2109 .. 2507) since the sheet was stored
# newest-first (descending) from row 16 to row 62.
$periods += "2508"

# ---------------------------------------------------------------------------
# Step 2: make room for the new row. Insert a blank row right after the old
# last data row (62) and before the signature block (old rows 67/68), which
# shifts the signature block down to rows 68/69.
# ---------------------------------------------------------------------------
$newRow = $lastDataRow + 1
$ws.Range("B$newRow`:J$newRow").Insert(-4121) | Out-Null   # xlShiftDown

# Give the brand-new row the special "closing row" formatting that the old
# last row (62) used to have.
$ws.Range("B$lastDataRow`:J$lastDataRow").Copy() | Out-Null
$ws.Range("B$newRow`:J$newRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Copy the repeating (non period) values into the new row.
$ws.Range("B$newRow").Value = $ws.Range("B$lastDataRow").Value()
$ws.Range("C$newRow").Value = $ws.Range("C$lastDataRow").Value()
$ws.Range("D$newRow").Value = $ws.Range("D$lastDataRow").Value()
$ws.Range("F$newRow").Value = $ws.Range("F$lastDataRow").Value()
$ws.Range("G$newRow").Value = $ws.Range("G$lastDataRow").Value()

# Now give the former last row (62) the regular interior-row formatting
# (copied from row 61, a standard data row).
$ws.Range("B61:J61").Copy() | Out-Null
$ws.Range("B$lastDataRow`:J$lastDataRow").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Step 3: write the ascending period list into column E for rows 16..63.
# ---------------------------------------------------------------------------
for ($i = 0; $i -lt $periods.Count; $i++) {
    $row = $firstDataRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# ---------------------------------------------------------------------------
# Step 4: update the summary fields.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 1744368   # Valor Mora
$ws.Range("F13").Value = 48        # Cant. Periodos
